# Generate Report for Handoff
#
# Inserts one new localization-status entry at the top of the data
# (just below the header row) and appends one new entry at the bottom,
# on all three worksheets (Overview, zh-cn, de-de), matching the
# "Generate Report for Handoff" CI report-generation commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Data describing the four data rows in their FINAL, post-edit order.
# Row 1 ("top") and Row 4 ("bottom") are the two brand-new entries;
# the middle two already existed in the workbook before the edit.
# ---------------------------------------------------------------------
$entries = @(
    @{ Guid = "1bf5a701-461b-4935-99d9-10695c2df1a9"; Hash = "5688db79eb84183e8b372b85cadd613690697902" },
    @{ Guid = "45ae3022-45ce-4049-9587-029da67ef300"; Hash = "5eed52a0e18d2a11ea66c61ac2b3196a42b734e0" },
    @{ Guid = "4874d6b2-94b7-42fe-bd9f-0ab8f849d697"; Hash = "6914c7948b82b51773bc729bee1f31d166d26695" },
    @{ Guid = "5d998178-445e-4313-bbf4-ee52016e9326"; Hash = "4408a069e20465cd18081b0f34c276952f20f1f4" }
)

$status            = "Ready for handoff"
$overviewDate      = "2016-02-18 04:02:07"
$zhHandoffDateTime = "2016-03-18 04:02:04"
$deHandoffDateTime = "2016-03-18 04:02:07"
$handbackDateTime  = "0001-01-01 00:00:00"
$handoffReason     = "Include"

$mdUrlTemplate  = "https://github.com/OpenLocalizationTest/oltest/blob/95331d3b62c39c34ba4f4f91d5495b5ecabbb285/e2e/{0}.md"
$zhUrlTemplate  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2cdf90d450b785b92be3ceb31bf05873f3227173/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/{0}.{1}.zh-cn.xlf"
$deUrlTemplate  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c97a9f8459908ab015c8748f204dd337678f877a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/{0}.{1}.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"  (columns A..D : File Name | zh-cn | de-de | Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

for ($i = 0; $i -lt $entries.Count; $i++) {
    $row  = $i + 2
    $guid = $entries[$i].Guid
    $mdUrl = [string]::Format($mdUrlTemplate, $guid)

    $aCell = $wsOverview.Cells.Item($row, 1)
    $wsOverview.Hyperlinks.Add($aCell, $mdUrl, "", "", "$guid.md") | Out-Null

    $wsOverview.Cells.Item($row, 2).Value = $status
    $wsOverview.Cells.Item($row, 3).Value = $status
    $wsOverview.Cells.Item($row, 4).Value = $overviewDate
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# columns: A Source File Name | B File Extension | C Status | D Latest Handoff File
#          E Latest Handoff Datetime | H Latest Handback DateTime | I Handoff Reason
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

for ($i = 0; $i -lt $entries.Count; $i++) {
    $row  = $i + 2
    $guid = $entries[$i].Guid
    $hash = $entries[$i].Hash
    $mdUrl  = [string]::Format($mdUrlTemplate, $guid)
    $xlfUrl = [string]::Format($zhUrlTemplate, $guid, $hash)
    $xlfName = "$guid.$hash.zh-cn.xlf"

    $aCell = $wsZh.Cells.Item($row, 1)
    $wsZh.Hyperlinks.Add($aCell, $mdUrl, "", "", "$guid.md") | Out-Null

    $bCell = $wsZh.Cells.Item($row, 2)
    $wsZh.Hyperlinks.Add($bCell, $mdUrl, "", "", ".md") | Out-Null

    $wsZh.Cells.Item($row, 3).Value = $status

    $dCell = $wsZh.Cells.Item($row, 4)
    $wsZh.Hyperlinks.Add($dCell, $xlfUrl, "", "", $xlfName) | Out-Null

    $wsZh.Cells.Item($row, 5).Value = $zhHandoffDateTime
    $wsZh.Cells.Item($row, 8).Value = $handbackDateTime
    $wsZh.Cells.Item($row, 9).Value = $handoffReason
}

# ---------------------------------------------------------------------
# Sheet 3: "de-de"  (same layout as "zh-cn")
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

for ($i = 0; $i -lt $entries.Count; $i++) {
    $row  = $i + 2
    $guid = $entries[$i].Guid
    $hash = $entries[$i].Hash
    $mdUrl  = [string]::Format($mdUrlTemplate, $guid)
    $xlfUrl = [string]::Format($deUrlTemplate, $guid, $hash)
    $xlfName = "$guid.$hash.de-de.xlf"

    $aCell = $wsDe.Cells.Item($row, 1)
    $wsDe.Hyperlinks.Add($aCell, $mdUrl, "", "", "$guid.md") | Out-Null

    $bCell = $wsDe.Cells.Item($row, 2)
    $wsDe.Hyperlinks.Add($bCell, $mdUrl, "", "", ".md") | Out-Null

    $wsDe.Cells.Item($row, 3).Value = $status

    $dCell = $wsDe.Cells.Item($row, 4)
    $wsDe.Hyperlinks.Add($dCell, $xlfUrl, "", "", $xlfName) | Out-Null

    $wsDe.Cells.Item($row, 5).Value = $deHandoffDateTime
    $wsDe.Cells.Item($row, 8).Value = $handbackDateTime
    $wsDe.Cells.Item($row, 9).Value = $handoffReason
}

Write-Output "Generate Report for Handoff: done"
